$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing existing rows 14-54 down to 15-55.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14. Values mirror the former row 14
# (now shifted to row 15) except Volumen (J) and Fecha (D), which take
# on the new reading for the added week.
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 44592
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = 300000001
$ws.Cells.Item(14, 7).Value = "Rabanito"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 50
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 6000
$ws.Cells.Item(14, 13).Value = 6000
$ws.Cells.Item(14, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(14, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(14, 16).Value = 500
$ws.Cells.Item(14, 17).Value = 12
$ws.Cells.Item(14, 18).Value = "Hortaliza"
